$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the existing header style
# (bold font + border + centered/top alignment) from E1 so the new column
# matches the other header cells exactly.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps captured when each panel entry was processed.
$timestamps = @(
    "2021-10-05 10:52:04.115876",
    "2021-10-05 10:52:04.115887",
    "2021-10-05 10:52:04.115891",
    "2021-10-05 10:52:04.115893",
    "2021-10-05 10:52:04.115896",
    "2021-10-05 10:52:04.115899",
    "2021-10-05 10:52:04.115901",
    "2021-10-05 10:52:04.115904",
    "2021-10-05 10:52:04.115907",
    "2021-10-05 10:52:04.115909",
    "2021-10-05 10:52:04.115912",
    "2021-10-05 10:52:04.115914",
    "2021-10-05 10:52:04.115917",
    "2021-10-05 10:52:04.115920",
    "2021-10-05 10:52:04.115922",
    "2021-10-05 10:52:04.115925",
    "2021-10-05 10:52:04.115927",
    "2021-10-05 10:52:04.115930",
    "2021-10-05 10:52:04.115933",
    "2021-10-05 10:52:04.115935"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Output "time_taken column added"
